$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Delete row 14 (SAL SOLUBLE) -- this shifts the old TOTAL row (15) up to row 14,
# preserving its style, and shrinks the used range to A1:F14.
$ws.Rows.Item(14).Delete()

# Update column widths (D, E, F) to the new target widths.
# NOTE: the ColumnWidth COM property reads/writes a value 5/6 of a character
# narrower than the raw OOXML <col width> attribute (Calibri 11 padding), so
# subtract 5/6 from the desired stored width to land exactly on target.
$ws.Columns.Item(4).ColumnWidth = 13 - 5/6
$ws.Columns.Item(5).ColumnWidth = 22 - 5/6
$ws.Columns.Item(6).ColumnWidth = 28 - 5/6

# Row 2 (240X120 PORCELANATO)
$ws.Range("D2").Value = 259.2
$ws.Range("E2").Value = 2633.00588040374
$ws.Range("F2").Value = 0.08962017598962102

# Row 3 (240X80 PORCELANATO)
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 20387.4774217135
$ws.Range("F3").Value = 0

# Row 4 (FREGADEROS DE COCINA)
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 782.417163948959
$ws.Range("F4").Value = 0

# Row 5 (GRIFERIAS) -- unchanged

# Row 6 (INODOROS)
$ws.Range("D6").Value = 321.3
$ws.Range("E6").Value = 2004.76694516821
$ws.Range("F6").Value = 0.1381301602980155

# Row 7 (LAVABOS)
$ws.Range("D7").Value = 390.37
$ws.Range("E7").Value = 496.341016287574
$ws.Range("F7").Value = 0.4402448969613309

# Row 8 (NO RESURTIBLES)
$ws.Range("D8").Value = 83.5
$ws.Range("E8").Value = 365.30162917203
$ws.Range("F8").Value = 0.1860510180278192

# Row 9 (OTROS)
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0

# Row 10 (PANELES DECORATIVOS)
$ws.Range("D10").Value = 56.32
$ws.Range("E10").Value = 2660.43588474074
$ws.Range("F10").Value = 0.0207306075294927

# Row 11 (PIEDRA SINTERIZADA)
$ws.Range("D11").Value = 648.83
$ws.Range("E11").Value = 18924.2302492497
$ws.Range("F11").Value = 0.03314913415365755

# Row 12 (PORCELANATO)
$ws.Range("C12").Value = 48624.06
$ws.Range("D12").Value = -43.78
$ws.Range("E12").Value = 48667.84
$ws.Range("F12").Value = -0.0009003773029236967

# Row 13 (PUERTAS DE SEGURIDAD)
$ws.Range("D13").Value = -124.24
$ws.Range("E13").Value = 1234.67665120341
$ws.Range("F13").Value = -0.1118839150935425

# Row 14 (was row 15, TOTAL) -- recalculated totals after the SAL SOLUBLE row removal
$ws.Range("C14").Value = 99897.99284188784
$ws.Range("D14").Value = 1591.5
$ws.Range("E14").Value = 98306.49284188786
$ws.Range("F14").Value = 0.01593125101641356
